# daily auto push: 2026-02-01 22:41 UTC
#
# Inserts a new data row right before the current row 740 (shifting the
# existing rows 740:781 down to 741:782) and fills the new row with the
# day's entry: 2026/02/02 (Mon), hour 5, ranking 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 740:781 down by one to make room for the new entry.
$ws.Rows("740:740").Insert()

# Column A (date) and B (weekday) are stored as plain text in this sheet,
# so force text format before assigning the value to stop "2026/02/02"
# from being auto-converted into a date serial, then drop back to the
# unstyled "Normal" cell style so the new cells match the rest of the
# sheet (which carries no explicit style).
$ws.Range("A740").NumberFormat = "@"
$ws.Range("A740").Value = "2026/02/02"
$ws.Range("A740").Style = "Normal"

$ws.Range("B740").NumberFormat = "@"
$ws.Range("B740").Value = "月"
$ws.Range("B740").Style = "Normal"

# Column C (time) and D (ranking) are plain numbers.
$ws.Range("C740").Value = 5
$ws.Range("D740").Value = 22
